$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bring column S's formatting across into the new column T (rows 3-40,
#    the bordered/bodied part of the table) so the new "2021" column looks
#    like the existing year columns.
$ws.Range("S3:S40").Copy()
$ws.Range("T3:T40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Fill in the 2021 figures (column T) - header year + data rows.
$ws.Range("T4").Value = 2021
$ws.Range("T6").Value = 1466
$ws.Range("T8").Value = 76
$ws.Range("T9").Value = 15
$ws.Range("T10").Value = 1
$ws.Range("T11").Value = 188
$ws.Range("T12").Value = 22
$ws.Range("T13").Value = 15
$ws.Range("T14").Value = "-"
$ws.Range("T15").Value = "-"
$ws.Range("T16").Value = 112
$ws.Range("T17").Value = "-"
$ws.Range("T18").Value = 6
$ws.Range("T19").Value = "-"
$ws.Range("T20").Value = 29
$ws.Range("T21").Value = 1002
$ws.Range("T22").Value = "-"
$ws.Range("T24").Value = 1029
$ws.Range("T26").Value = 51
$ws.Range("T27").Value = 4
$ws.Range("T28").Value = "-"
$ws.Range("T29").Value = 127
$ws.Range("T30").Value = 14
$ws.Range("T31").Value = 12
$ws.Range("T32").Value = "-"
$ws.Range("T33").Value = "-"
$ws.Range("T34").Value = 70
$ws.Range("T35").Value = "-"
$ws.Range("T36").Value = 3
$ws.Range("T37").Value = "-"
$ws.Range("T38").Value = 16
$ws.Range("T39").Value = 732
$ws.Range("T40").Value = "-"


# 3) The two subtotal rows ("Всего детей...") use a bold right-aligned
#    style in column S; replicate that on the new column for those rows.
$ws.Range("T6").HorizontalAlignment = -4152
$ws.Range("T24").HorizontalAlignment = -4152

# 4) Match the saved selection state of the authored workbook.
$ws.Range("T3").Select() | Out-Null
